# Update of Excel Modules Files
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "ElementName3" header - the column itself had no other data.
$ws.Cells.Item(1, 6).ClearContents()

# Fill in the new data columns (B..E) for each module row, plus the code
# column stays the same (GIL51..GIL56).
$data = @(
    @("GIL51", "pede. Suspendisse dui.",       "EL Haddad",        "Nullam feugiat placerat",   "varius et, euismod"),
    @("GIL52", "a nunc. In",                   "Badir",             "sodales nisi magna",         "elementum sem, vitae"),
    @("GIL53", "amet metus. Aliquam",          "Ezzine",            "Cras vulputate velit",       "scelerisque neque sed"),
    @("GIL54", "quam vel sapien",              "El Alami Hassoun",  "Nunc mauris elit,",          "libero et tristique"),
    @("GIL55", "feugiat nec, diam.",           "Lazaar",            "pellentesque. Sed dictum.",  "ridiculus mus. Proin"),
    @("GIL54-2", "nonummy. Fusce fermentum",   "El Haddad",         "neque pellentesque massa",   "Mauris eu turpis."),
    @("GIL55-2", "a, arcu. Sed",               "EL Haddad",         "sit amet risus.",            "Nulla facilisi. Sed"),
    @("GIL56", "Suspendisse eleifend. Cras",   "El Alami Hassoun",  "velit dui, semper",          "ligula elit, pretium")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

# B2 keeps the formatting that came along with the pasted text (explicit
# black font color rather than the theme color).
$ws.Cells.Item(2, 2).Font.Color = 0

# Autofit the columns to the new content.
$ws.Columns("A:E").AutoFit() | Out-Null

# Leave the selection where the user ended up after entering the data.
$ws.Range("H10").Select()
